# Replaced any doc strings where "hour" is used with "time" to show
# generalisation of functions for any reporting interval.
#
# Visible workbook-level effects of that change (as captured by the
# regenerated test fixture):
#   1. The "Criterion 1, Air Speed 0.1" sheet tab is now listed after
#      "Criterion 2, Air Speed 0.1" (tabs rotate: 1,3,2 -> 3,2,1).
#   2. The readme index table's column order changes from
#      (index, JobNo, Date, sheet_name, Author) to
#      (index, JobNo, Author, Date, sheet_name) and the run Date is
#      bumped from 2022-03-16 to 2022-03-24.

$wb = $excel.ActiveWorkbook

# --- 1. Move "Criterion 1, Air Speed 0.1" to the end of the tab list ---
$sheets = $wb.Worksheets
$moving = $sheets.Item("Criterion 1, Air Speed 0.1")
$lastSheet = $sheets.Item($sheets.Count)
$moving.Move([System.Reflection.Missing]::Value, $lastSheet)

# --- 2. Rebuild the readme table (Table1) with the new column order/date ---
$readme = $wb.Worksheets.Item("readme")

# Header row: B stays JobNo, C/D/E swap from Date/sheet_name/Author
# to Author/Date/sheet_name.
$readme.Range("C1").Value = "Author"
$readme.Range("D1").Value = "Date"
$readme.Range("E1").Value = "sheet_name"

# Data rows keep the same "index" (column A) per row, but the sheet_name
# reflects the new tab order, and Author/Date move into C/D.
$readme.Range("C2").Value = "jovyan"
$readme.Range("D2").Value = "'20220324"
$readme.Range("E2").Value = "Criteria Failing, Air Speed 0.1"

$readme.Range("C3").Value = "jovyan"
$readme.Range("D3").Value = "'20220324"
$readme.Range("E3").Value = "Criterion 3, Air Speed 0.1"

$readme.Range("C4").Value = "jovyan"
$readme.Range("D4").Value = "'20220324"
$readme.Range("E4").Value = "Criterion 2, Air Speed 0.1"

$readme.Range("C5").Value = "jovyan"
$readme.Range("D5").Value = "'20220324"
$readme.Range("E5").Value = "Criterion 1, Air Speed 0.1"
